# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows before the previous last row of data,
# which shifts the old row 360 down to row 364.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(360).Insert()
$ws.Rows.Item(360).Insert()
$ws.Rows.Item(360).Insert()
$ws.Rows.Item(360).Insert()

# Common values shared by all rows in this block.
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"

# Row 360: Lapins, Primera
$r = 360
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44890
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Lapins"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 330
$ws.Cells.Item($r, 14).Value = 14400
$ws.Cells.Item($r, 15).Value = 14400
$ws.Cells.Item($r, 16).Value = 14400
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 800
$ws.Cells.Item($r, 20).Value = 18

# Row 361: Lapins, Segunda
$r = 361
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44890
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Lapins"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 12600
$ws.Cells.Item($r, 15).Value = 12600
$ws.Cells.Item($r, 16).Value = 12600
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 700
$ws.Cells.Item($r, 20).Value = 18

# Row 362: Royal Dawn, Primera
$r = 362
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44890
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Royal Dawn"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = 10

# Row 363: Santina, Primera
$r = 363
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44890
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Santina"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 350
$ws.Cells.Item($r, 14).Value = 8000
$ws.Cells.Item($r, 15).Value = 8000
$ws.Cells.Item($r, 16).Value = 8000
$ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($r, 19).Value = 800
$ws.Cells.Item($r, 20).Value = 10
